$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original text formatting
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '50.937.01'
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").Value = '2.887.22'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '368.38'
$ws.Range("E5").Value = '  +4.06%  '
$ws.Range("D6").Value = '101.87'
$ws.Range("E6").Value = '  -4.94%  '
$ws.Range("E7").Value = '  -4.72%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -5.38%  '
$ws.Range("D10").Value = '36.50'
$ws.Range("E10").Value = '  -4.36%  '
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").Value = '0.0830'
$ws.Range("E12").Value = '  -3.11%  '
$ws.Range("D13").Value = '18.26'
$ws.Range("E13").Value = '  -4.64%  '
$ws.Range("D14").Value = '3.341.51'
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("D15").Value = '7.33'
$ws.Range("E15").Value = '  -4.53%  '
$ws.Range("D16").Value = '2.884.47'
$ws.Range("E16").Value = '  -3.01%  '
$ws.Range("D17").Value = '0.926'
$ws.Range("E17").Value = '  -4.52%  '
$ws.Range("D18").Value = '50.846.34'
$ws.Range("E18").Value = '  -2.13%  '
$ws.Range("D19").Value = '3.24'
$ws.Range("E19").Value = '  -5.52%  '
$ws.Range("D20").Value = '7.17'
$ws.Range("E20").Value = '  -4.33%  '
$ws.Range("D21").Value = '12.71'
$ws.Range("E21").Value = '  -6.43%  '
$ws.Range("D22").Value = '0.0₃0937'
$ws.Range("E22").Value = '  -3.58%  '
$ws.Range("D23").Value = '67.84'
$ws.Range("E23").Value = '  -2.87%  '
$ws.Range("D24").Value = '257.91'
$ws.Range("E24").Value = '  -3.10%  '
$ws.Range("D25").Value = '2.66'
$ws.Range("E25").Value = '  -3.04%  '
$ws.Range("E26").Value = '  -3.15%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").Value = '25.48'
$ws.Range("E28").Value = '  -4.99%  '
$ws.Range("D29").Value = '6.95'
$ws.Range("E29").Value = '  -8.56%  '
$ws.Range("D30").Value = '0.101'
$ws.Range("E30").Value = '  -3.04%  '
$ws.Range("D31").Value = '9.82'
$ws.Range("E31").Value = '  -4.70%  '
$ws.Range("D32").Value = '6.03'
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("D33").Value = '2.12'
$ws.Range("E33").Value = '  -2.43%  '
$ws.Range("D34").Value = '34.24'
$ws.Range("E34").Value = '  -6.66%  '
$ws.Range("D35").Value = '50.98'
$ws.Range("E35").Value = '  -1.99%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Value = '0.0414'
$ws.Range("E37").Value = '  -4.10%  '
$ws.Range("E38").Value = '  -4.65%  '
$ws.Range("D39").Value = '2.61'
$ws.Range("E39").Value = '  -2.44%  '
$ws.Range("D40").Value = '16.88'
$ws.Range("E40").Value = '  -5.75%  '
$ws.Range("D41").Value = '1.83'
$ws.Range("E41").Value = '  -7.56%  '
$ws.Range("E42").Value = '  -5.07%  '
$ws.Range("D43").Value = '21.81'
$ws.Range("E43").Value = '  -4.46%  '
$ws.Range("D44").Value = '118.48'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("E45").Value = '  -3.31%  '
$ws.Range("D46").Value = '2.008.58'
$ws.Range("E46").Value = '  -5.05%  '
$ws.Range("E47").Value = '  -6.59%  '
$ws.Range("D48").Value = '3.13'
$ws.Range("E48").Value = '  -6.78%  '
$ws.Range("D49").Value = '3.183.45'
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("D50").Value = '0.233'
$ws.Range("E50").Value = '  -3.37%  '
$ws.Range("E51").Value = '  -10.50%  '
